$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-unused label cells (text removed, formatting/style preserved)
$ws.Range("A2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("A3").Value = ""

# C6 picks up the same "horizontal left" format that A2 uses (the old,
# now-redundant border-less left-aligned style), its text stays unchanged
$ws.Range("A2").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null

# Move the active selection to F6
$ws.Range("F6").Select() | Out-Null
